# Apply the "17.12.2020" daily update to the Slovakia Covid daily stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrections to already-existing rows (AgTests / AgPosit revisions) ---
$ws.Range("H273").Value = 26052

$ws.Range("H274").Value = 27352
$ws.Range("I274").Value = 1262

$ws.Range("H275").Value = 27416
$ws.Range("I275").Value = 1186

$ws.Range("H276").Value = 12551
$ws.Range("I276").Value = 424

$ws.Range("H278").Value = 28741
$ws.Range("I278").Value = 1957

$ws.Range("H279").Value = 42606
$ws.Range("I279").Value = 3008

$ws.Range("H280").Value = 33854
$ws.Range("I280").Value = 2244

$ws.Range("H281").Value = 42544
$ws.Range("I281").Value = 3040

$ws.Range("H282").Value = 43293
$ws.Range("I282").Value = 2612

$ws.Range("H283").Value = 16537
$ws.Range("I283").Value = 960

$ws.Range("H284").Value = 1357
$ws.Range("I284").Value = 114

$ws.Range("H285").Value = 37513
$ws.Range("I285").Value = 3065

$ws.Range("H286").Value = 50718
$ws.Range("I286").Value = 3793

# --- New row for 2020-12-16 (serial date 44181) ---
$ws.Range("A287").Value = 44181
$ws.Range("A287").NumberFormat = $ws.Range("A286").NumberFormat
$ws.Range("B287").Value = 142133
$ws.Range("C287").Value = 102737
$ws.Range("D287").Value = 38018
$ws.Range("E287").Value = 14921
$ws.Range("F287").Value = 3045
$ws.Range("G287").Value = 1378
$ws.Range("H287").Value = 49741
$ws.Range("I287").Value = 3369

$wb.Save()
